# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns.
# D-column cells hold number-like text (e.g. "24.50", "1.000") that must
# stay text, not be reinterpreted as numeric values by Excel, so each is
# written as Text (NumberFormat "@") and the style is reset to Normal
# immediately after so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.396.87"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.08%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.849.20"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.18%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "240.65"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.6298"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.07662"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +1.71%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.2942"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "24.50"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.04%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07750"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.839.68"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("E14").Value = "  +8.94%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.6798"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.56%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "83.56"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.92%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.094.58"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -7.50%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "6.138"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "29.428.02"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "229.25"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.73%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "12.46"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.448"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("E24").Value = "  +0.02%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "157.22"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.1391"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "17.67"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.468"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  +3.63%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.05648"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.71%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "4.111"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.048"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.85%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.851"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.7098"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.586"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.777"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.229.61"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.08%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.01800"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  +3.91%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9114"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.000"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.003.81"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.09%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "101.46"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.48%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "66.16"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "7.162"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.71%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.4013"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "9.069"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.690"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "
